# Apply updated crypto price/volume data (commit: "Updated cryptos list on Tue Jun 27 09:14:07 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "30.375.51"
$ws.Cells.Item(2, 5).Value = "  +0.05%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.875.67"
$ws.Cells.Item(3, 5).Value = "  -0.79%  "

$ws.Cells.Item(4, 5).Value = "  -0.02%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "238.27"
$ws.Cells.Item(5, 5).Value = "  +0.13%  "

$ws.Cells.Item(6, 5).Value = "  +0.02%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4804"
$ws.Cells.Item(7, 5).Value = "  -0.41%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2814"
$ws.Cells.Item(8, 5).Value = "  -3.00%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06507"
$ws.Cells.Item(9, 5).Value = "  -1.42%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.870.47"
$ws.Cells.Item(10, 5).Value = "  -1.11%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07479"
$ws.Cells.Item(11, 5).Value = "  +1.15%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "16.56"
$ws.Cells.Item(12, 5).Value = "  -2.03%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "5.070"
$ws.Cells.Item(13, 5).Value = "  -2.04%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "88.19"
$ws.Cells.Item(14, 5).Value = "  +0.62%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.6598"
$ws.Cells.Item(15, 5).Value = "  -0.35%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "30.349.59"
$ws.Cells.Item(16, 5).Value = "  -0.03%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "13.28"
$ws.Cells.Item(17, 5).Value = "  -1.24%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "1.000"
$ws.Cells.Item(18, 5).Value = "  +0.04%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000007570"
$ws.Cells.Item(19, 5).Value = "  -2.53%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "2.118.03"
$ws.Cells.Item(20, 5).Value = "  -1.01%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "5.297"
$ws.Cells.Item(21, 5).Value = "  -3.40%  "

$ws.Cells.Item(22, 5).Value = "  -0.01%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "220.72"
$ws.Cells.Item(23, 5).Value = "  +14.38%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "6.180"
$ws.Cells.Item(24, 5).Value = "  -0.12%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.330"
$ws.Cells.Item(25, 5).Value = "  -1.12%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "166.65"
$ws.Cells.Item(26, 5).Value = "  +0.98%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "18.50"
$ws.Cells.Item(27, 5).Value = "  +1.14%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.962"
$ws.Cells.Item(28, 5).Value = "  -0.22%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.463"
$ws.Cells.Item(29, 5).Value = "  +0.94%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.09359"
$ws.Cells.Item(30, 5).Value = "  +2.21%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.306"
$ws.Cells.Item(31, 5).Value = "  +0.96%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.021"
$ws.Cells.Item(32, 5).Value = "  -0.82%  "

$ws.Cells.Item(33, 5).Value = "  -1.26%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.197"
$ws.Cells.Item(34, 5).Value = "  +4.65%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.7410"
$ws.Cells.Item(35, 5).Value = "  +0.97%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.710"
$ws.Cells.Item(36, 5).Value = "  -0.13%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.01818"
$ws.Cells.Item(37, 5).Value = "  +0.18%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.616"
$ws.Cells.Item(38, 5).Value = "  -1.12%  "

$ws.Cells.Item(39, 2).Value = "TrustWalletToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.9056"
$ws.Cells.Item(39, 5).Value = "  -1.50%  "

$ws.Cells.Item(40, 2).Value = "RenderToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.055"
$ws.Cells.Item(40, 5).Value = "  -1.45%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "106.61"
$ws.Cells.Item(41, 5).Value = "  -0.04%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "5.881"
$ws.Cells.Item(42, 5).Value = "  -0.20%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.4263"
$ws.Cells.Item(43, 5).Value = "  -1.81%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "7.391"
$ws.Cells.Item(45, 5).Value = "  -2.67%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "64.08"
$ws.Cells.Item(46, 5).Value = "  -1.25%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.1272"
$ws.Cells.Item(47, 5).Value = "  -4.73%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.473"
$ws.Cells.Item(48, 5).Value = "  -5.10%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "8.914"
$ws.Cells.Item(49, 5).Value = "  -1.05%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "33.61"
$ws.Cells.Item(50, 5).Value = "  -1.71%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.3881"
$ws.Cells.Item(51, 5).Value = "  +0.16%  "
